# Scheduled-runner profit-sheet refresh: updates recalculated crafting-profit
# figures (currentAveragePrice / LevePrice / Profit columns, H:N) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets to match freshly pulled
# market-board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 119.666664
$ws.Range("I5").Value = 92.5
$ws.Range("K5").Value = 92.5
$ws.Range("M5").Value = 22.5

# Row 43
$ws.Range("H43").Value = 36112736
$ws.Range("I43").Value = 216666670
$ws.Range("K43").Value = 216666670
$ws.Range("M43").Value = -216666601

# Row 135
$ws.Range("H135").Value = 350.7037
$ws.Range("I135").Value = 362.8
$ws.Range("J135").Value = 199.5
$ws.Range("K135").Value = 3265.2
$ws.Range("L135").Value = 1795.5
$ws.Range("M135").Value = -730.2000000000003
$ws.Range("N135").Value = -6865.5

# Row 141
$ws.Range("H141").Value = 1678.6666
$ws.Range("I141").Value = 1441.4286
$ws.Range("K141").Value = 4324.2858
$ws.Range("M141").Value = 855.7142000000003

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 266182.1
$ws.Range("I32").Value = 333863.7
$ws.Range("K32").Value = 333863.7
$ws.Range("M32").Value = -333576.7

# Row 45
$ws.Range("H45").Value = 29541.389
$ws.Range("I45").Value = 35057.535
$ws.Range("K45").Value = 35057.535
$ws.Range("M45").Value = -34680.535

# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 76923730
$ws.Range("J80").Value = 240.14285
$ws.Range("L80").Value = 240.14285
$ws.Range("N80").Value = -2236.14285

# Row 83
$ws.Range("H83").Value = 76923730
$ws.Range("J83").Value = 240.14285
$ws.Range("L83").Value = 1200.71425
$ws.Range("N83").Value = -11184.71425

# Row 97
$ws.Range("H97").Value = 26146.715
$ws.Range("I97").Value = 23837.834
$ws.Range("J97").Value = 40000
$ws.Range("K97").Value = 23837.834
$ws.Range("L97").Value = 40000
$ws.Range("M97").Value = -22846.834
$ws.Range("N97").Value = -41982

# Row 105
$ws.Range("H105").Value = 28411.334
$ws.Range("I105").Value = 23093.6
$ws.Range("J105").Value = 55000
$ws.Range("K105").Value = 23093.6
$ws.Range("L105").Value = 55000
$ws.Range("M105").Value = -21346.6
$ws.Range("N105").Value = -58494

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 10007500
$ws.Range("I4").Value = 15000
$ws.Range("K4").Value = 15000
$ws.Range("M4").Value = -14888

# Row 16
$ws.Range("H16").Value = 17859470
$ws.Range("I16").Value = 35715332
$ws.Range("J16").Value = 3608
$ws.Range("K16").Value = 35715332
$ws.Range("L16").Value = 3608
$ws.Range("M16").Value = -35715045
$ws.Range("N16").Value = -4182

# Row 31
$ws.Range("H31").Value = 3043.8125
$ws.Range("I31").Value = 3188.353
$ws.Range("J31").Value = 2880
$ws.Range("K31").Value = 3188.353
$ws.Range("L31").Value = 2880
$ws.Range("M31").Value = -2893.353
$ws.Range("N31").Value = -3470

# Row 34
$ws.Range("H34").Value = 3043.8125
$ws.Range("I34").Value = 3188.353
$ws.Range("J34").Value = 2880
$ws.Range("K34").Value = 3188.353
$ws.Range("L34").Value = 2880
$ws.Range("M34").Value = -2986.353
$ws.Range("N34").Value = -3284

# Row 58
$ws.Range("H58").Value = 3088.4546
$ws.Range("J58").Value = 3996.5715
$ws.Range("L58").Value = 3996.5715
$ws.Range("N58").Value = -4402.5715

# Row 86
$ws.Range("H86").Value = 12193.833
$ws.Range("J86").Value = 18273.889
$ws.Range("L86").Value = 18273.889
$ws.Range("N86").Value = -20519.889

# Row 89
$ws.Range("H89").Value = 12193.833
$ws.Range("J89").Value = 18273.889
$ws.Range("L89").Value = 91369.44499999999
$ws.Range("N89").Value = -102601.445

# Row 113
$ws.Range("H113").Value = 17859470
$ws.Range("I113").Value = 35715332
$ws.Range("J113").Value = 3608
$ws.Range("K113").Value = 35715332
$ws.Range("L113").Value = 3608
$ws.Range("M113").Value = -35713162
$ws.Range("N113").Value = -7948

# Row 134
$ws.Range("H134").Value = 2562.75
$ws.Range("I134").Value = 2589.4285
$ws.Range("K134").Value = 7768.2855
$ws.Range("M134").Value = -5233.2855

# Row 136
$ws.Range("H136").Value = 3088.4546
$ws.Range("J136").Value = 3996.5715
$ws.Range("L136").Value = 11989.7145
$ws.Range("N136").Value = -17089.7145

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1149.25
$ws.Range("I5").Value = 904.3333
$ws.Range("J5").Value = 1464.1428
$ws.Range("K5").Value = 2712.9999
$ws.Range("L5").Value = 4392.428400000001
$ws.Range("M5").Value = -2600.9999
$ws.Range("N5").Value = -4616.428400000001

# Row 18
$ws.Range("H18").Value = 413.33334
$ws.Range("I18").Value = 413.33334
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1240.00002
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1071.00002
$ws.Range("N18").ClearContents()

# Row 122
$ws.Range("H122").Value = 6061984
$ws.Range("I122").Value = 9524284
$ws.Range("K122").Value = 85718556
$ws.Range("M122").Value = -85716106

# Row 135
$ws.Range("H135").Value = 1149.25
$ws.Range("I135").Value = 904.3333
$ws.Range("J135").Value = 1464.1428
$ws.Range("K135").Value = 8138.9997
$ws.Range("L135").Value = 13177.2852
$ws.Range("M135").Value = -5603.9997
$ws.Range("N135").Value = -18247.2852

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 897.5
$ws.Range("I3").Value = 897.5
$ws.Range("K3").Value = 897.5
$ws.Range("M3").Value = -781.5

# Row 10
$ws.Range("H10").Value = 2323896.8
$ws.Range("J10").Value = 2323896.8
$ws.Range("L10").Value = 2323896.8
$ws.Range("N10").Value = -2324234.8

# Row 20
$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10490

# Row 45
$ws.Range("H45").Value = 75000
$ws.Range("J45").Value = 75000
$ws.Range("L45").Value = 75000
$ws.Range("N45").Value = -76118

# Row 97
$ws.Range("H97").Value = 1119.3
$ws.Range("I97").Value = 574.5
$ws.Range("K97").Value = 574.5
$ws.Range("M97").Value = -78.5

# Row 122
$ws.Range("H122").Value = 3011.2083
$ws.Range("I122").Value = 2814.1
$ws.Range("K122").Value = 8442.299999999999
$ws.Range("M122").Value = -5992.299999999999

# Row 126
$ws.Range("H126").Value = 2305.6
$ws.Range("I126").Value = 2244
$ws.Range("J126").Value = 2398
$ws.Range("K126").Value = 6732
$ws.Range("L126").Value = 7194
$ws.Range("M126").Value = -4262
$ws.Range("N126").Value = -12134

# Row 128
$ws.Range("H128").Value = 84999.5
$ws.Range("J128").Value = 84999.5
$ws.Range("L128").Value = 84999.5
$ws.Range("N128").Value = -94959.5

# Row 130
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 1019.4
$ws.Range("I55").Value = 1196.1111
$ws.Range("K55").Value = 1196.1111
$ws.Range("M55").Value = -1023.1111

# Row 136
$ws.Range("H136").Value = 52975.05
$ws.Range("I136").Value = 128204.25
$ws.Range("J136").Value = 2822.25
$ws.Range("K136").Value = 384612.75
$ws.Range("L136").Value = 8466.75
$ws.Range("M136").Value = -382062.75
$ws.Range("N136").Value = -13566.75

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3284.8572
$ws.Range("I62").Value = 3124.75
$ws.Range("J62").Value = 3498.3333
$ws.Range("K62").Value = 3124.75
$ws.Range("L62").Value = 3498.3333
$ws.Range("M62").Value = -2500.75
$ws.Range("N62").Value = -4746.3333

# Row 65
$ws.Range("H65").Value = 3284.8572
$ws.Range("I65").Value = 3124.75
$ws.Range("J65").Value = 3498.3333
$ws.Range("K65").Value = 15623.75
$ws.Range("L65").Value = 17491.6665
$ws.Range("M65").Value = -12503.75
$ws.Range("N65").Value = -23731.6665
